# OCR FLOW COMPLETED / MODIFICATION IN EDUCATION FLOW
# Appends new education-sprint rows (13-18) below the existing data on the
# "AMSIN" sheet, re-styles row 12 to match the rest of the table body, and
# refreshes the run-time stamp in B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMSIN")

# --- Row 12: pick up the common body style (s="5") on the text/number cells
#     (column B already carries the date-time style) and refresh the run
#     timestamp to the corrected value.
$ws.Range("A12").Style = "Normal"
$ws.Range("C12:G12").Style = "Normal"

$ws.Cells.Item(12, 1).Value = "2022-09-01"
$ws.Cells.Item(12, 2).Value2 = 44805.42799498843
$ws.Cells.Item(12, 3).Value = "vbn1278"
$ws.Cells.Item(12, 4).Value2 = 60
$ws.Cells.Item(12, 5).Value2 = 60
$ws.Cells.Item(12, 6).Value2 = 0
$ws.Cells.Item(12, 7).Value2 = 2.74

# --- New rows of registration-history data (rows 13-17 use the common body
#     style; row 18 -- the newest entry -- keeps the default/no-style
#     formatting, exactly like row 12 did before this edit).
$newRows = @(
    @{ Row = 13; Date = "2022-09-02"; Time = 44806.67115123843; Name = "edu223";    Total = 60; Pass = 53; Fail = 7;  Taken = 1.76 },
    @{ Row = 14; Date = "2022-09-02"; Time = 44806.67537111111; Name = "wew344";    Total = 60; Pass = 28; Fail = 32; Taken = 0.39 },
    @{ Row = 15; Date = "2022-09-02"; Time = 44806.67871225694; Name = "ssd33";     Total = 60; Pass = 32; Fail = 28; Taken = 0.7 },
    @{ Row = 16; Date = "2022-09-02"; Time = 44806.68225763889; Name = "vinod234";  Total = 60; Pass = 21; Fail = 39; Taken = 3.37 },
    @{ Row = 17; Date = "2022-09-02"; Time = 44806.68892793982; Name = "edu654";    Total = 60; Pass = 58; Fail = 2;  Taken = 1.7 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Style = "Normal"
    $ws.Range("C$row`:G$row").Style = "Normal"

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value2 = $r.Time
    $ws.Cells.Item($row, 3).Value = $r.Name
    $ws.Cells.Item($row, 4).Value2 = $r.Total
    $ws.Cells.Item($row, 5).Value2 = $r.Pass
    $ws.Cells.Item($row, 6).Value2 = $r.Fail
    $ws.Cells.Item($row, 7).Value2 = $r.Taken
}

# Row 18 -- default (no explicit cell style) formatting, like legacy row 12.
$ws.Cells.Item(18, 1).Value = "2022-09-02"
$ws.Cells.Item(18, 2).Value2 = 44806.70102041208
$ws.Cells.Item(18, 3).Value = "e166"
$ws.Cells.Item(18, 4).Value2 = 60
$ws.Cells.Item(18, 5).Value2 = 60
$ws.Cells.Item(18, 6).Value2 = 0
$ws.Cells.Item(18, 7).Value2 = 1.16

# The B column keeps its date-time number format on every row.
$ws.Range("B13:B18").NumberFormat = $ws.Range("B12").NumberFormat
